$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")
$ws.Activate()
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1, 5).Value = "End time"
$ws.Range("E1").Select() | Out-Null
